$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pt_max column (E2:E12) values from 50 to 70
$ws.Range("E2:E12").Value = 70

# Update the sheet selection to match the edited state
$ws.Range("E2:E12").Select()
